$wb = $excel.ActiveWorkbook

# --- Overview sheet: update shared status text (zh-cn status in E3, de-de status in F3) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handback transform failed"
$wsOverview.Range("F3").Value = "Handback transform failed"

# --- zh-cn sheet: same status text (Status column, C3), new Error Detail (P3), widen column P ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handback transform failed"
$wsZhCn.Range("P3").Value = "Handback file name: a34iq1ks.tnh is different with handoff file name: da8fc9c6-9eee-4c0b-9c71-35c776b47446.0ab5a1cd0be8403a4f43b57f7ac0575b9672c226.zh-cn."
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

# --- de-de sheet: same status text (Status column, C3), new Error Detail (P3), widen column P ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handback transform failed"
$wsDeDe.Range("P3").Value = "Handback file name: a34iq1ks.tnh is different with handoff file name: da8fc9c6-9eee-4c0b-9c71-35c776b47446.0ab5a1cd0be8403a4f43b57f7ac0575b9672c226.de-de."
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
